$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = ""
$ws.Range("A1").Style = "Normal"
$ws.Range("A1").Value = 'questions = [
    {
        "title": "You are an SEO outreach consultant for a wine delivery service. A new competitor has recently emerged and is quickly gaining traffic and rankings for a lot of similar keywords to those your business targets.Which action should you take?",
        "ques_type": 2,
        "options": [
            "Use a backlink analysis tool to identify the domains linking to the competitor''s website.",
            "Look for patterns in the anchor text of the competitor''s backlinks and replicate them for your website.",
            "Reach out to the competitor''s PR team and ask for a list of their current backlinks.",
            "Scrape the competitor''s website and analyze the internal links to identify potential link-building opportunities."
        ],
        "score": "Use a backlink analysis tool to identify the domains linking to the competitor''s website."
    },
    {
        "title": "The health and diet supplement company you work for as an SEO outreach consultant is launching a new weight-loss product. Because this is such a competitive market, you want to build an effective link-building campaign to quickly increase the product\u2019s visibility.Which action should you take?",
        "ques_type": 2,
        "options": [
            "Add backlinks to the product''s webpage in several third-party comments sections.",
            "Participate in link exchange schemes with other weight-loss product companies.",
            "Guest blog on relevant health and fitness websites.",
            "Create new websites that backlink to the new product''s webpage."
        ],
        "score": "Guest blog on relevant health and fitness websites."
    },
    {
        "title": "The real estate company you are carrying out SEO outreach for is looking to expand its link-building efforts. You have found a third-party website for backlinking but need to identify whether it would be a good fit. Which action should you take?",
        "ques_type": 2,
        "options": [
            "Use the search query \"link:example.com\" to find if the website has recently linked to your competitors.",
            "Use the search query \"intitle:keyword\" to find all websites that have recently published related content.",
            "Use the search query \"inurl:brandname\" to find all websites that have recently mentioned your company''s brand name.",
            "Use the search query \"site:example.com\" + \"keyword\" to find pages that might be relevant to your company."
        ],
        "score": "Use the search query \"site:example.com\" + \"keyword\" to find pages that might be relevant to your company."
    },
    {
        "title": "You are an SEO outreach specialist for a time-tracking app. The link-building campaign you are running has been live for a few weeks, and you have limited resources to evaluate how effective it has been so far.Which of the following actions should you take?",
        "ques_type": 2,
        "options": [
            "Monitor social media engagement and customer reviews.",
            "Analyze the increase in organic traffic, organic rankings, and authority of backlinks.",
            "Conduct surveys and collect user feedback.",
            "Review only the number of backlinks acquired during the campaign."
        ],
        "score": "Analyze the increase in organic traffic, organic rankings, and authority of backlinks."
    }
]'
$ws.Rows(1).AutoFit()
